$d = $word.ActiveDocument

# Target paragraph: " * לעשות דף מנהל ששם יהיה את הטבלה של כל המשתמשים"
# (the note about building an admin page whose table lists all users).
# Mark it as done/resolved by striking through the whole line, matching
# the sibling notes below it that already use strikethrough.
$p = $d.Paragraphs.Item(25)
$p.Range.Font.StrikeThrough = $true
